# Update the numeric results for the "case with 380 kV" re-run (pl_mw.xlsx, Sheet1).
# Columns B,C,D,F,G,H,J,L,M,N,O for rows 2-25 get new values; A (index) and the
# always-zero columns E,I,K are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B", "C", "D", "F", "G", "H", "J", "L", "M", "N", "O")

$data = @{
    2 = @('0.977845674447849', '0.2734622245895366', '0.2214358252317083', '1.253651696097307', '0.6567907814089722', '0.7701826252495749', '0.2229225274845356', '0.3280651365297373', '0.26721371665883', '1.465218122115189', '2.845898825768586')
    3 = @('0.9024912028600625', '0.2671379641161451', '0.2214603815762857', '1.257586859865661', '0.6576568856440375', '0.774394196528057', '0.2246085862268781', '0.3250958909660824', '0.2547157732088507', '1.475217591408736', '2.856278980734061')
    4 = @('0.8563620134113137', '0.2632431640283812', '0.2215425151552211', '1.260636837119307', '0.6586048167407554', '0.7773058055078863', '0.2257003439741544', '0.3233802269507819', '0.247112011517153', '1.481856619471635', '2.864206653534254')
    5 = @('0.8376003736272821', '0.2616532101689586', '0.2215929034100057', '1.262039229308023', '0.6590957436239648', '0.7785742903925978', '0.2261594737772035', '0.3227082335410643', '0.2440313022478122', '1.484687795938697', '2.867828224232611')
    6 = @('0.8344872586376084', '0.261389035244548', '0.2216022938152022', '1.262281732938028', '0.6591835813631377', '0.7787898752366331', '0.2262365719431894', '0.3225982934041838', '0.2435208411968901', '1.485165508562183', '2.868453202194402')
    7 = @('0.8561088381061666', '0.2632217324628385', '0.2215431261325662', '1.260655104275664', '0.658611013886329', '0.7773225806914255', '0.2257064783178555', '0.3233710540949986', '0.2470703912104355', '1.481894292460247', '2.864253912031074')
    8 = @('0.9518355205230762', '0.2712841287283396', '0.2214304082736405', '1.254877080356572', '0.6570030400611842', '0.7715672284826809', '0.2234921635973206', '0.3270191132409863', '0.2628900535591754', '1.468562403113154', '2.849155416834549')
    9 = @('1.140600316245809', '0.2869961461368575', '0.2217390975996878', '1.248570159641218', '0.6571527631309664', '0.7628617569975802', '0.2195973594517786', '0.3350204899314377', '0.2944574940046536', '1.446373447701042', '2.831874731227686')
    10 = @('1.279861539576586', '0.2984733168020171', '0.2222859152077845', '1.246993961458642', '0.659279082691242', '0.7580349270876354', '0.2170072800279771', '0.3414097630552817', '0.317970760300625', '1.432472997163615', '2.826691157924529')
    11 = @('1.343327907992091', '0.3036788381735107', '0.2226035641045954', '1.246939827195476', '0.6606848879796843', '0.7561789285696676', '0.2158876361464247', '0.3444260727443691', '0.3287348596375423', '1.426668982260317', '2.825963991660416')
    12 = @('1.367376279382825', '0.3056476708004823', '0.2227337030995713', '1.247014542180466', '0.6612803191542156', '0.7555248915075623', '0.2154720620747159', '0.3455839390338014', '0.3328204435475257', '1.424545698697457', '2.825923061851199')
    13 = @('1.362196389221822', '0.3052237560380604', '0.2227052380370296', '1.246994218008496', '0.6611492764026963', '0.7556635812312464', '0.2155611895575169', '0.3453338776629096', '0.3319401242326308', '1.424999671014547', '2.82592145175056')
    14 = @('1.34530608986006', '0.3038408639190209', '0.2226140735804023', '1.246944066422991', '0.6607326100522073', '0.7561241431170203', '0.2158532781400154', '0.3445210181097451', '0.3290707953247747', '1.426492804570259', '2.825955926460949')
    15 = @('1.334962204113708', '0.3029934872921842', '0.2225595141328824', '1.246925743648795', '0.6604856058375503', '0.7564126024070532', '0.2160332857098459', '0.3440251529831215', '0.3273144708120981', '1.427417100006956', '2.826007570670953')
    16 = @('1.275716054610598', '0.2981327984429214', '0.2222665377954485', '1.247010829031325', '0.6591960363745244', '0.7581630519254929', '0.2170816291832354', '0.3412148381154765', '0.3172686400289209', '1.432862744419175', '2.826771490004063')
    17 = @('1.239398959795324', '0.2951468416997614', '0.2221044135568349', '1.247232730794373', '0.6585172624640023', '0.7593238669962545', '0.2177397509100277', '0.3395188277016103', '0.3111230018917581', '1.436336411757154', '2.827657774666932')
    18 = @('1.218521317369209', '0.2934279486281923', '0.2220176523692743', '1.247422763299355', '0.6581681250644209', '0.7600235200217043', '0.2181238013800845', '0.3385536699807403', '0.3075945923899184', '1.438383265437317', '2.828321052299685')
    19 = @('1.21145444293694', '0.2928457167487437', '0.2219893925950771', '1.247497825654357', '0.6580570015141944', '0.7602659056348529', '0.2182547821622208', '0.3382286645817629', '0.3064010416173915', '1.439084695720574', '2.828571994760836')
    20 = @('1.243263850769381', '0.2954648530671591', '0.2221210008367152', '1.24720265153524', '0.6585852470868332', '0.7591969866210349', '0.2176691219281097', '0.3396983014512074', '0.3117765554341219', '1.435961574625473', '2.827547541031521')
    21 = @('1.350266784970188', '0.3042471188364573', '0.2226405838442602', '1.246956213929536', '0.6608532828228988', '0.755987541328409', '0.2157672564980712', '0.3447593508478661', '0.3299133325814907', '1.426052212106256', '2.825939438663312')
    22 = @('1.420286278860601', '0.3099728613971138', '0.2230375622318661', '1.247350084814826', '0.6627032823370342', '0.7541743338944826', '0.2145732949017085', '0.3481582282709468', '0.3418217276450548', '1.420010477045537', '2.826254860592201')
    23 = @('1.382908137620859', '0.3069182551411984', '0.2228204534667739', '1.24708912801124', '0.6616822503153656', '0.7551160806526553', '0.2152060545918655', '0.3463358836746693', '0.3354610617128202', '1.423195334900399', '2.825961515917015')
    24 = @('1.241516528641853', '0.2953210870529261', '0.2221134816449961', '1.247216055808032', '0.6585543832210874', '0.7592542486677729', '0.2177010355734024', '0.3396171305178513', '0.31148106890452', '1.436130883372556', '2.827596898760135')
    25 = @('1.089428907602496', '0.2827568296088856', '0.2215991793793251', '1.249739047788395', '0.6567583672940316', '0.7649409667448737', '0.2206032437394174', '0.3327657934384547', '0.2858605113060335', '1.451953721668161', '2.835230139778417')
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Count; $i++) {
        $cellRef = "{0}{1}" -f $columns[$i], $row
        $ws.Range($cellRef).Value = $values[$i]
    }
}
